$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.6011661417734331
$ws.Cells.Item(2, 3).Value = 0.1741463629819222
$ws.Cells.Item(2, 5).Value = 0.125274054624299
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.002482616468633616
$ws.Cells.Item(2, 11).Value = 0.2822869884435306
$ws.Cells.Item(2, 12).Value = 0.1913268763976106
$ws.Cells.Item(2, 13).Value = 0.1693692603555874
$ws.Cells.Item(2, 14).Value = 2.111026505566134
$ws.Cells.Item(2, 15).Value = 3.91108259546607

# Row 3
$ws.Cells.Item(3, 2).Value = 0.5691298085733649
$ws.Cells.Item(3, 3).Value = 0.1740889422873835
$ws.Cells.Item(3, 5).Value = 0.1256336016998034
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.002484820388283449
$ws.Cells.Item(3, 11).Value = 0.253726555037602
$ws.Cells.Item(3, 12).Value = 0.1887399731340338
$ws.Cells.Item(3, 13).Value = 0.1632560345410958
$ws.Cells.Item(3, 14).Value = 2.130938582997276
$ws.Cells.Item(3, 15).Value = 3.939404973576814

# Row 4
$ws.Cells.Item(4, 2).Value = 0.5496775708365362
$ws.Cells.Item(4, 3).Value = 0.1740673321471427
$ws.Cells.Item(4, 5).Value = 0.1259028539400422
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.002486246453175626
$ws.Cells.Item(4, 11).Value = 0.2362229106123834
$ws.Cells.Item(4, 12).Value = 0.1872412534128074
$ws.Cells.Item(4, 13).Value = 0.1595745942517439
$ws.Cells.Item(4, 14).Value = 2.143799970077141
$ws.Cells.Item(4, 15).Value = 3.958588135631658

# Row 5
$ws.Cells.Item(5, 2).Value = 0.541806085311606
$ws.Cells.Item(5, 3).Value = 0.1740619776840724
$ws.Cells.Item(5, 5).Value = 0.1260247931725402
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.002486845957662681
$ws.Cells.Item(5, 11).Value = 0.2290986494674172
$ws.Cells.Item(5, 12).Value = 0.1866531285908124
$ws.Cells.Item(5, 13).Value = 0.158092630725541
$ws.Cells.Item(5, 14).Value = 2.14920090930137
$ws.Cells.Item(5, 15).Value = 3.966856523878704

# Row 6
$ws.Cells.Item(6, 2).Value = 0.5405023981365389
$ws.Cells.Item(6, 3).Value = 0.1740612976630587
$ws.Cells.Item(6, 5).Value = 0.1260457796521823
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.002486946615939073
$ws.Cells.Item(6, 11).Value = 0.2279162042125904
$ws.Cells.Item(6, 12).Value = 0.1865568389757328
$ws.Cells.Item(6, 13).Value = 0.1578476578505423
$ws.Cells.Item(6, 14).Value = 2.150107385962827
$ws.Cells.Item(6, 15).Value = 3.968256734404278

# Row 7
$ws.Cells.Item(7, 2).Value = 0.5495711878938323
$ws.Cells.Item(7, 3).Value = 0.1740672459326369
$ws.Cells.Item(7, 5).Value = 0.1259044489581296
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.002486254463742954
$ws.Cells.Item(7, 11).Value = 0.236126794886502
$ws.Cells.Item(7, 12).Value = 0.1872332300896602
$ws.Cells.Item(7, 13).Value = 0.1595545339048137
$ws.Cells.Item(7, 14).Value = 2.143872161834913
$ws.Cells.Item(7, 15).Value = 3.95869781927648

# Row 8
$ws.Cells.Item(8, 2).Value = 0.5900750702529933
$ws.Cells.Item(8, 3).Value = 0.1741237440969456
$ws.Cells.Item(8, 5).Value = 0.1253879767339612
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.00248336129233587
$ws.Cells.Item(8, 11).Value = 0.2724328352850307
$ws.Cells.Item(8, 12).Value = 0.1904163477414826
$ws.Cells.Item(8, 13).Value = 0.1672465255761857
$ws.Cells.Item(8, 14).Value = 2.117760379327644
$ws.Cells.Item(8, 15).Value = 3.920476112189149

# Row 9
$ws.Cells.Item(9, 2).Value = 0.671212767496371
$ws.Cells.Item(9, 3).Value = 0.1743419611043322
$ws.Cells.Item(9, 5).Value = 0.1247589345773861
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.002478263396254747
$ws.Cells.Item(9, 11).Value = 0.343872227034808
$ws.Cells.Item(9, 12).Value = 0.1973670864061745
$ws.Cells.Item(9, 13).Value = 0.1828981642628804
$ws.Cells.Item(9, 14).Value = 2.071592110719408
$ws.Cells.Item(9, 15).Value = 3.859742973842543

# Row 10
$ws.Cells.Item(10, 2).Value = 0.7318442784064985
$ws.Cells.Item(10, 3).Value = 0.17456668488974
$ws.Cells.Item(10, 5).Value = 0.1245295104141597
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.002474865478927584
$ws.Cells.Item(10, 11).Value = 0.3964919494622734
$ws.Cells.Item(10, 12).Value = 0.2029029776374927
$ws.Cells.Item(10, 13).Value = 0.19473882590286
$ws.Cells.Item(10, 14).Value = 2.040736478911785
$ws.Cells.Item(10, 15).Value = 3.823781544200585

# Row 11
$ws.Cells.Item(11, 2).Value = 0.7596441069940454
$ws.Cells.Item(11, 3).Value = 0.1746826757015683
$ws.Cells.Item(11, 5).Value = 0.1244754245238084
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.002473394414087672
$ws.Cells.Item(11, 11).Value = 0.420456035855949
$ws.Cells.Item(11, 12).Value = 0.2055140359681786
$ws.Cells.Item(11, 13).Value = 0.2001986714289714
$ws.Cells.Item(11, 14).Value = 2.027363449722238
$ws.Cells.Item(11, 15).Value = 3.809300196985276

# Row 12
$ws.Cells.Item(12, 2).Value = 0.7702020459667835
$ws.Cells.Item(12, 3).Value = 0.1747285578909725
$ws.Cells.Item(12, 5).Value = 0.1244621514952691
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.002472848042208749
$ws.Cells.Item(12, 11).Value = 0.4295341453808419
$ws.Cells.Item(12, 12).Value = 0.2065160506116541
$ws.Cells.Item(12, 13).Value = 0.2022766313711699
$ws.Cells.Item(12, 14).Value = 2.022394755022832
$ws.Cells.Item(12, 15).Value = 3.804086349014199

# Row 13
$ws.Cells.Item(13, 2).Value = 0.7679268475567937
$ws.Cells.Item(13, 3).Value = 0.1747185895311745
$ws.Cells.Item(13, 5).Value = 0.1244646898484447
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.002472965238591849
$ws.Cells.Item(13, 11).Value = 0.4275788669732492
$ws.Cells.Item(13, 12).Value = 0.2062996601167129
$ws.Cells.Item(13, 13).Value = 0.2018286434906997
$ws.Cells.Item(13, 14).Value = 2.023460610321145
$ws.Cells.Item(13, 15).Value = 3.805197240871337

# Row 14
$ws.Cells.Item(14, 2).Value = 0.7605121024170103
$ws.Cells.Item(14, 3).Value = 0.1746864113108728
$ws.Cells.Item(14, 5).Value = 0.1244741882089251
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.002473349249634143
$ws.Cells.Item(14, 11).Value = 0.4212028313301346
$ws.Cells.Item(14, 12).Value = 0.2055962069215411
$ws.Cells.Item(14, 13).Value = 0.2003694180566526
$ws.Cells.Item(14, 14).Value = 2.02695276067336
$ws.Cells.Item(14, 15).Value = 3.808865841382215

# Row 15
$ws.Cells.Item(15, 2).Value = 0.7559743422095266
$ws.Cells.Item(15, 3).Value = 0.174666955743362
$ws.Cells.Item(15, 5).Value = 0.1244809442831425
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.002473585858261261
$ws.Cells.Item(15, 11).Value = 0.4172977581872033
$ws.Cells.Item(15, 12).Value = 0.2051670466202467
$ws.Cells.Item(15, 13).Value = 0.1994769554246076
$ws.Cells.Item(15, 14).Value = 2.0291042248567
$ws.Cells.Item(15, 15).Value = 3.811148114673131

# Row 16
$ws.Cells.Item(16, 2).Value = 0.7300318315461709
$ws.Cells.Item(16, 3).Value = 0.174559379840332
$ws.Cells.Item(16, 5).Value = 0.1245340549314413
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.002474963114812065
$ws.Cells.Item(16, 11).Value = 0.3949263455624248
$ws.Cells.Item(16, 12).Value = 0.2027341989247162
$ws.Cells.Item(16, 13).Value = 0.1943834794490371
$ws.Cells.Item(16, 14).Value = 2.041623785415656
$ws.Cells.Item(16, 15).Value = 3.824765710863147

# Row 17
$ws.Cells.Item(17, 2).Value = 0.7141723875971024
$ws.Cells.Item(17, 3).Value = 0.1744968961643707
$ws.Cells.Item(17, 5).Value = 0.1245794988438043
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.002475827106323801
$ws.Cells.Item(17, 11).Value = 0.3812088190616407
$ws.Cells.Item(17, 12).Value = 0.2012654292558267
$ws.Cells.Item(17, 13).Value = 0.1912775217938005
$ws.Cells.Item(17, 14).Value = 2.049474024180082
$ws.Cells.Item(17, 15).Value = 3.833600524918722

# Row 18
$ws.Cells.Item(18, 2).Value = 0.7050710357364096
$ws.Cells.Item(18, 3).Value = 0.1744622538236129
$ws.Cells.Item(18, 5).Value = 0.1246103709239748
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.002476331081849296
$ws.Cells.Item(18, 11).Value = 0.3733214404412877
$ws.Cells.Item(18, 12).Value = 0.2004293666577297
$ws.Cells.Item(18, 13).Value = 0.1894979801362098
$ws.Cells.Item(18, 14).Value = 2.054051712061773
$ws.Cells.Item(18, 15).Value = 3.83885881815317

# Row 19
$ws.Cells.Item(19, 2).Value = 0.7019930285575526
$ws.Cells.Item(19, 3).Value = 0.1744507478306403
$ws.Cells.Item(19, 5).Value = 0.1246216374133553
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.002476502928443169
$ws.Cells.Item(19, 11).Value = 0.3706513669167748
$ws.Cells.Item(19, 12).Value = 0.20014779271159
$ws.Cells.Item(19, 13).Value = 0.1888966508996148
$ws.Cells.Item(19, 14).Value = 2.05561236263302
$ws.Cells.Item(19, 15).Value = 3.840669544095562

# Row 20
$ws.Cells.Item(20, 2).Value = 0.7158585270529727
$ws.Cells.Item(20, 3).Value = 0.174503413609294
$ws.Cells.Item(20, 5).Value = 0.124574171476791
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.002475734405695873
$ws.Cells.Item(20, 11).Value = 0.3826688095668942
$ws.Cells.Item(20, 12).Value = 0.2014208788069709
$ws.Cells.Item(20, 13).Value = 0.1916074407204107
$ws.Cells.Item(20, 14).Value = 2.048631890811148
$ws.Cells.Item(20, 15).Value = 3.832641752378578

# Row 21
$ws.Cells.Item(21, 2).Value = 0.7626891634071455
$ws.Cells.Item(21, 3).Value = 0.174695809826666
$ws.Cells.Item(21, 5).Value = 0.1244712028600183
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.002473236166433891
$ws.Cells.Item(21, 11).Value = 0.4230755381985887
$ws.Cells.Item(21, 12).Value = 0.2058024688018492
$ws.Cells.Item(21, 13).Value = 0.2007977457732082
$ws.Cells.Item(21, 14).Value = 2.025924442895534
$ws.Cells.Item(21, 15).Value = 3.807780959501002

# Row 22
$ws.Cells.Item(22, 2).Value = 0.7934746139481774
$ws.Cells.Item(22, 3).Value = 0.1748329611097645
$ws.Cells.Item(22, 5).Value = 0.1244459111196008
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.002471665702508368
$ws.Cells.Item(22, 11).Value = 0.4495034468257302
$ws.Cells.Item(22, 12).Value = 0.2087433664862033
$ws.Cells.Item(22, 13).Value = 0.2068649038236785
$ws.Cells.Item(22, 14).Value = 2.011639765984436
$ws.Cells.Item(22, 15).Value = 3.79310627853269

# Row 23
$ws.Cells.Item(23, 2).Value = 0.7770276787282455
$ws.Cells.Item(23, 3).Value = 0.1747587235491537
$ws.Cells.Item(23, 5).Value = 0.1244555737175439
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.002472498205894568
$ws.Cells.Item(23, 11).Value = 0.4353967199913029
$ws.Cells.Item(23, 12).Value = 0.2071667079535331
$ws.Cells.Item(23, 13).Value = 0.203621230562689
$ws.Cells.Item(23, 14).Value = 2.019212899102673
$ws.Cells.Item(23, 15).Value = 3.800794504254725

# Row 24
$ws.Cells.Item(24, 2).Value = 0.7150961721050351
$ws.Cells.Item(24, 3).Value = 0.1745004630821754
$ws.Cells.Item(24, 5).Value = 0.124576565195472
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.002475776292918235
$ws.Cells.Item(24, 11).Value = 0.3820087509078007
$ws.Cells.Item(24, 12).Value = 0.2013505740422374
$ws.Cells.Item(24, 13).Value = 0.1914582653172587
$ws.Cells.Item(24, 14).Value = 2.049012418450745
$ws.Cells.Item(24, 15).Value = 3.83307465573597

# Row 25
$ws.Cells.Item(25, 2).Value = 0.649082174276117
$ws.Cells.Item(25, 3).Value = 0.1742715401467159
$ws.Cells.Item(25, 5).Value = 0.1248881489826328
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.002479581242531691
$ws.Cells.Item(25, 11).Value = 0.3245214917721455
$ws.Cells.Item(25, 12).Value = 0.1954111532822651
$ws.Cells.Item(25, 13).Value = 0.1786036902453034
$ws.Cells.Item(25, 14).Value = 2.083543294011156
$ws.Cells.Item(25, 15).Value = 3.87465157670502
